$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two "no survey data yet" placeholder rows near the top
#     (old rows 4 & 5, reference_number 3 & 4) and near the middle
#     (old rows 19 & 20, reference_number 18 & 19). Delete bottom-up so
#     the remaining row numbers don't shift under us while we work.
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# After the four deletions above, the sheet now ends at row 20
# (reference_number 23) exactly as in the target workbook.

# --- Append the new Nupku sites. Cell values are written in the same
#     left-to-right, first-seen order used when the rows were originally
#     authored so new shared-string entries land in the matching slots.

# Row 21: Hartley Creek / 197542_us
$ws.Cells.Item(21,1).Value = 35
$ws.Cells.Item(21,2).Value = "Hartley Creek"
$ws.Cells.Item(21,3).Value = "197542_us"
$ws.Cells.Item(21,4).Value = 725
$ws.Cells.Item(21,8).Value = "high"
$ws.Cells.Item(21,10).Value = "high"
$ws.Cells.Item(21,13).Value = 7200

# Row 22: Hartley Creek / 197542_ds
$ws.Cells.Item(22,1).Value = 36
$ws.Cells.Item(22,2).Value = "Hartley Creek"
$ws.Cells.Item(22,3).Value = "197542_ds"
$ws.Cells.Item(22,4).Value = 900
$ws.Cells.Item(22,8).Value = "medium"

# Row 23: Hartley Creek / 197582_ds
$ws.Cells.Item(23,1).Value = 37
$ws.Cells.Item(23,2).Value = "Hartley Creek"
$ws.Cells.Item(23,3).Value = "197582_ds"
$ws.Cells.Item(23,4).Value = 540
$ws.Cells.Item(23,8).Value = "medium"

# Row 24: Unnamed Tributary to Morrissey Creek / 50181_us
$ws.Cells.Item(24,1).Value = 38
$ws.Cells.Item(24,2).Value = "Unnamed Tributary to Morrissey Creek"
$ws.Cells.Item(24,3).Value = "50181_us"
$ws.Cells.Item(24,4).Value = 515
$ws.Cells.Item(24,8).Value = "medium"
$ws.Cells.Item(24,10).Value = "high"
$ws.Cells.Item(24,13).Value = 515

# Row 25: Unnamed Tributary to Morrissey Creek / 50181_ds
$ws.Cells.Item(25,1).Value = 39
$ws.Cells.Item(25,2).Value = "Unnamed Tributary to Morrissey Creek"
$ws.Cells.Item(25,3).Value = "50181_ds"
$ws.Cells.Item(25,4).Value = 200
$ws.Cells.Item(25,8).Value = "medium"

# Row 26: Stove Creek / 50152_us
$ws.Cells.Item(26,1).Value = 40
$ws.Cells.Item(26,2).Value = "Stove Creek"
$ws.Cells.Item(26,3).Value = "50152_us"
$ws.Cells.Item(26,4).Value = 675
$ws.Cells.Item(26,8).Value = "high"
$ws.Cells.Item(26,10).Value = "high"
$ws.Cells.Item(26,13).Value = 2700

# Row 28: Weigart Creek / 197534_us
$ws.Cells.Item(28,1).Value = 42
$ws.Cells.Item(28,2).Value = "Weigart Creek"
$ws.Cells.Item(28,3).Value = "197534_us"
$ws.Cells.Item(28,4).Value = 1100
$ws.Cells.Item(28,8).Value = "high"
$ws.Cells.Item(28,10).Value = "high"
$ws.Cells.Item(28,13).Value = 11600

# Row 29: Weigart Creek / 197534_ds
$ws.Cells.Item(29,1).Value = 43
$ws.Cells.Item(29,2).Value = "Weigart Creek"
$ws.Cells.Item(29,3).Value = "197534_ds"
$ws.Cells.Item(29,4).Value = 675
$ws.Cells.Item(29,8).Value = "medium"

# Row 30: Unnamed Tributary to Michel Creek / 50261_us
$ws.Cells.Item(30,1).Value = 44
$ws.Cells.Item(30,2).Value = "Unnamed Tributary to Michel Creek"
$ws.Cells.Item(30,3).Value = "50261_us"
$ws.Cells.Item(30,4).Value = 220
$ws.Cells.Item(30,8).Value = "medium"
$ws.Cells.Item(30,13).Value = 0

# Row 31: Unnamed Tributary to Michel Creek / 50261_ds
$ws.Cells.Item(31,1).Value = 45
$ws.Cells.Item(31,2).Value = "Unnamed Tributary to Michel Creek"
$ws.Cells.Item(31,3).Value = "50261_ds"
$ws.Cells.Item(31,4).Value = 210
$ws.Cells.Item(31,8).Value = "medium"

# Finish row 30's priority column and row 27 (Stove Creek / 50152_ds),
# written last to match the order these values were entered originally.
$ws.Cells.Item(30,10).Value = "no fix"
$ws.Cells.Item(27,1).Value = 41
$ws.Cells.Item(27,2).Value = "Stove Creek"
$ws.Cells.Item(27,3).Value = "50152_ds"
$ws.Cells.Item(27,4).Value = 140
$ws.Cells.Item(27,8).Value = "high"

# --- View tidy-up: scroll back to column A and move the active
#     selection to N24, matching the saved sheet view.
$ws.Range("A1").Select()
$ws.Range("N24").Select()

"done"
